$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices, volumes, and a few reordered rows)
$updates = @{
    'D2' = '98.222.12'
    'E2' = '  +4.37%  '
    'D3' = '3.368.48'
    'E3' = '  +9.67%  '
    'E4' = '  +0.01%  '
    'D5' = '254.75'
    'E5' = '  +8.11%  '
    'D6' = '623.20'
    'E6' = '  +2.22%  '
    'E7' = '  +8.47%  '
    'E8' = '  +2.17%  '
    'E9' = '  -0.01%  '
    'D10' = '3.367.68'
    'E10' = '  +9.63%  '
    'D11' = '0.818'
    'E11' = '  +0.82%  '
    'E12' = '  +1.33%  '
    'D13' = '97.974.47'
    'E13' = '  +4.20%  '
    'D14' = '35.91'
    'E14' = '  +6.17%  '
    'E15' = '  +2.70%  '
    'D16' = '3.990.73'
    'E16' = '  +9.26%  '
    'E17' = '  +3.40%  '
    'D18' = '3.375.39'
    'E18' = '  +9.97%  '
    'D19' = '3.66'
    'E19' = '  +2.98%  '
    'D20' = '14.85'
    'E20' = '  +3.58%  '
    'D21' = '482.26'
    'E21' = '  +8.60%  '
    'D22' = '5.91'
    'E22' = '  +3.21%  '
    'E23' = '  +10.04%  '
    'D24' = '9.24'
    'E24' = '  +4.86%  '
    'D25' = '5.80'
    'E25' = '  +5.54%  '
    'D26' = '88.16'
    'E26' = '  +4.28%  '
    'D27' = '12.05'
    'E27' = '  +1.07%  '
    'D28' = '3.544.91'
    'E28' = '  +9.34%  '
    'B30' = 'Stellar'
    'C30' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D30' = '0.252'
    'E30' = '  +1.02%  '
    'B31' = 'Cronos'
    'C31' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D31' = '0.186'
    'E31' = '  +3.79%  '
    'E32' = '  +2.49%  '
    'D33' = '0.999'
    'E33' = '  +0.03%  '
    'E34' = '  +3.62%  '
    'D35' = '27.40'
    'E35' = '  +7.76%  '
    'D36' = '527.52'
    'E36' = '  +10.17%  '
    'D37' = '0.153'
    'E37' = '  +1.66%  '
    'D38' = '7.36'
    'E38' = '  -2.30%  '
    'D39' = '1.96'
    'E39' = '  +3.99%  '
    'B40' = 'MantraDAO'
    'C40' = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
    'D40' = '3.93'
    'E40' = '  +4.58%  '
    'B41' = 'WhiteBITCoin'
    'C41' = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
    'D41' = '24.81'
    'E41' = '  +3.21%  '
    'B42' = 'PolygonEcosystemToken'
    'C42' = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
    'D42' = '0.449'
    'E42' = '  +2.99%  '
    'E43' = '  +1.97%  '
    'B44' = 'ARBITRUM'
    'C44' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D44' = '0.791'
    'E44' = '  +17.87%  '
    'B45' = 'dogwifhat'
    'C45' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D45' = '3.24'
    'E45' = '  +5.14%  '
    'D47' = '161.20'
    'E47' = '  -0.22%  '
    'D48' = '1.94'
    'E48' = '  +6.86%  '
    'B49' = 'ImmutableX'
    'C49' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D49' = '1.37'
    'E49' = '  +6.94%  '
    'B50' = 'OKB'
    'C50' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D50' = '45.51'
    'E50' = '  +4.34%  '
    'B51' = 'Filecoin'
    'C51' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D51' = '4.55'
    'E51' = '  +7.09%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
